$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price / 1h-volume-change figures (and two name/link swaps)
# as scraped on the latest GitHub Actions run.
# Note: some Price values (column D) look like plain numbers (e.g. "0.999"),
# so they are prefixed with a leading apostrophe to force Excel to store them
# as text, matching the original workbook's text-formatted Price column.
$ws.Range("D2").Value = '69.107.36'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '3.769.16'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''628.44'
$ws.Range("E5").Value = '  +3.31%  '
$ws.Range("D6").Value = '''165.97'
$ws.Range("E6").Value = '  +1.03%  '
$ws.Range("D7").Value = '3.767.20'
$ws.Range("E7").Value = '  -1.10%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").Value = '''6.78'
$ws.Range("E12").Value = '  -3.02%  '
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").Value = '''34.88'
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").Value = '4.404.39'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '3.753.96'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").Value = '69.104.70'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("E18").Value = '  -3.30%  '
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("E20").Value = '  -1.32%  '
$ws.Range("D21").Value = '''463.19'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = '''9.51'
$ws.Range("E22").Value = '  -1.19%  '
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").Value = '''82.94'
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("D26").Value = '''11.97'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("D28").Value = '''10.04'
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = '3.918.95'
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''2.27'
$ws.Range("E31").Value = '  +1.68%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.67'
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("E33").Value = '  -2.29%  '
$ws.Range("D34").Value = '''28.51'
$ws.Range("E34").Value = '  -2.34%  '
$ws.Range("D35").Value = '''0.169'
$ws.Range("E35").Value = '  +13.65%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("B37").Value = 'RenzoRestakedETH'
$ws.Range("C37").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D37").Value = '3.723.63'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '''8.98'
$ws.Range("E38").Value = '  -0.85%  '
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("E40").Value = '  +2.47%  '
$ws.Range("D41").Value = '''5.78'
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '''157.64'
$ws.Range("E45").Value = '  +3.07%  '
$ws.Range("D46").Value = '''1.42'
$ws.Range("E46").Value = '  +1.66%  '
$ws.Range("E47").Value = '  +3.79%  '
$ws.Range("D48").Value = '''43.07'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = '''46.59'
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("E51").Value = '  -0.54%  '
